# Auto-generated: update market-price derived columns (H-N) for specific Leve rows
# across multiple job-class sheets, matching a scheduled price-refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 294.22223
$ws.Cells.Item(39, 9).Value = 120.2
$ws.Cells.Item(39, 10).Value = 511.75
$ws.Cells.Item(39, 11).Value = 360.6
$ws.Cells.Item(39, 12).Value = 1535.25
$ws.Cells.Item(39, 13).Value = -64.60000000000002
$ws.Cells.Item(39, 14).Value = -2127.25
$ws.Cells.Item(103, 8).Value = 1445.8
$ws.Cells.Item(103, 9).Value = 1101.1428
$ws.Cells.Item(103, 10).Value = 2250
$ws.Cells.Item(103, 11).Value = 3303.4284
$ws.Cells.Item(103, 12).Value = 6750
$ws.Cells.Item(103, 13).Value = -2717.4284
$ws.Cells.Item(103, 14).Value = -7922
$ws.Cells.Item(117, 8).Value = 48654.5
$ws.Cells.Item(117, 10).Value = 48654.5
$ws.Cells.Item(117, 12).Value = 48654.5
$ws.Cells.Item(117, 14).Value = -57832.5
$ws.Cells.Item(133, 8).Value = 45888.777
$ws.Cells.Item(133, 10).Value = 45888.777
$ws.Cells.Item(133, 12).Value = 45888.777
$ws.Cells.Item(133, 14).Value = -56008.777
$ws.Cells.Item(138, 8).Value = 1701.7177
$ws.Cells.Item(138, 9).Value = 990.9737
$ws.Cells.Item(138, 11).Value = 2972.9211
$ws.Cells.Item(138, 13).Value = 2167.0789

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11216.2
$ws.Cells.Item(32, 9).Value = 10225.49
$ws.Cells.Item(32, 10).Value = 19307
$ws.Cells.Item(32, 11).Value = 10225.49
$ws.Cells.Item(32, 12).Value = 19307
$ws.Cells.Item(32, 13).Value = -9938.49
$ws.Cells.Item(32, 14).Value = -19881
$ws.Cells.Item(80, 8).Value = 38314.11
$ws.Cells.Item(80, 10).Value = 38314.11
$ws.Cells.Item(80, 12).Value = 38314.11
$ws.Cells.Item(80, 14).Value = -40310.11
$ws.Cells.Item(83, 8).Value = 38314.11
$ws.Cells.Item(83, 10).Value = 38314.11
$ws.Cells.Item(83, 12).Value = 114942.33
$ws.Cells.Item(83, 14).Value = -124926.33
$ws.Cells.Item(114, 8).Value = 45941.75
$ws.Cells.Item(114, 10).Value = 45941.75
$ws.Cells.Item(114, 12).Value = 45941.75
$ws.Cells.Item(114, 14).Value = -54619.75
$ws.Cells.Item(121, 8).Value = 39981
$ws.Cells.Item(121, 10).Value = 39981
$ws.Cells.Item(121, 12).Value = 39981
$ws.Cells.Item(121, 14).Value = -43475
$ws.Cells.Item(131, 8).Value = 44346
$ws.Cells.Item(131, 10).Value = 44346
$ws.Cells.Item(131, 12).Value = 44346
$ws.Cells.Item(131, 14).Value = -54426
$ws.Cells.Item(132, 8).Value = 2302.8928
$ws.Cells.Item(132, 9).Value = 1331.1
$ws.Cells.Item(132, 11).Value = 3993.3
$ws.Cells.Item(132, 13).Value = -1463.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(57, 8).Value = 55169.5
$ws.Cells.Item(57, 10).Value = 55169.5
$ws.Cells.Item(57, 12).Value = 55169.5
$ws.Cells.Item(57, 14).Value = -56609.5
$ws.Cells.Item(86, 8).Value = 1901
$ws.Cells.Item(86, 9).Value = 2302
$ws.Cells.Item(86, 10).Value = 1299.5
$ws.Cells.Item(86, 11).Value = 2302
$ws.Cells.Item(86, 12).Value = 1299.5
$ws.Cells.Item(86, 13).Value = -1179
$ws.Cells.Item(86, 14).Value = -3545.5
$ws.Cells.Item(89, 8).Value = 1901
$ws.Cells.Item(89, 9).Value = 2302
$ws.Cells.Item(89, 10).Value = 1299.5
$ws.Cells.Item(89, 11).Value = 11510
$ws.Cells.Item(89, 12).Value = 6497.5
$ws.Cells.Item(89, 13).Value = -5894
$ws.Cells.Item(89, 14).Value = -17729.5
$ws.Cells.Item(133, 8).Value = 48566.668
$ws.Cells.Item(133, 10).Value = 48566.668
$ws.Cells.Item(133, 12).Value = 48566.668
$ws.Cells.Item(133, 14).Value = -58686.668
$ws.Cells.Item(136, 8).Value = 55169.5
$ws.Cells.Item(136, 10).Value = 55169.5
$ws.Cells.Item(136, 12).Value = 55169.5
$ws.Cells.Item(136, 14).Value = -65369.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(52, 8).Value = 65799.664
$ws.Cells.Item(52, 10).Value = 65799.664
$ws.Cells.Item(52, 12).Value = 65799.664
$ws.Cells.Item(52, 14).Value = -66387.664
$ws.Cells.Item(137, 8).Value = 35184.617
$ws.Cells.Item(137, 10).Value = 35184.617
$ws.Cells.Item(137, 12).Value = 35184.617
$ws.Cells.Item(137, 14).Value = -45384.617
$ws.Cells.Item(139, 8).Value = 59239.8
$ws.Cells.Item(139, 10).Value = 63549.75
$ws.Cells.Item(139, 12).Value = 63549.75
$ws.Cells.Item(139, 14).Value = -73829.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4310.8335
$ws.Cells.Item(80, 9).Value = 4449.5
$ws.Cells.Item(80, 10).Value = 4137.5
$ws.Cells.Item(80, 11).Value = 4449.5
$ws.Cells.Item(80, 12).Value = 4137.5
$ws.Cells.Item(80, 13).Value = -3451.5
$ws.Cells.Item(80, 14).Value = -6133.5
$ws.Cells.Item(83, 8).Value = 4310.8335
$ws.Cells.Item(83, 9).Value = 4449.5
$ws.Cells.Item(83, 10).Value = 4137.5
$ws.Cells.Item(83, 11).Value = 22247.5
$ws.Cells.Item(83, 12).Value = 20687.5
$ws.Cells.Item(83, 13).Value = -17255.5
$ws.Cells.Item(83, 14).Value = -30671.5
$ws.Cells.Item(102, 8).Value = 1514.4073
$ws.Cells.Item(102, 9).Value = 1607.6086
$ws.Cells.Item(102, 10).Value = 978.5
$ws.Cells.Item(102, 11).Value = 1607.6086
$ws.Cells.Item(102, 12).Value = 978.5
$ws.Cells.Item(102, 13).Value = 14.39139999999998
$ws.Cells.Item(102, 14).Value = -4222.5
$ws.Cells.Item(110, 8).Value = 47008.668
$ws.Cells.Item(110, 10).Value = 47008.668
$ws.Cells.Item(110, 12).Value = 47008.668
$ws.Cells.Item(110, 14).Value = -55188.668
$ws.Cells.Item(119, 8).Value = 48442
$ws.Cells.Item(119, 10).Value = 48442
$ws.Cells.Item(119, 12).Value = 48442
$ws.Cells.Item(119, 14).Value = -58118

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2723.75
$ws.Cells.Item(68, 9).Value = 1860
$ws.Cells.Item(68, 10).Value = 3116.3635
$ws.Cells.Item(68, 11).Value = 1860
$ws.Cells.Item(68, 12).Value = 3116.3635
$ws.Cells.Item(68, 13).Value = -1111
$ws.Cells.Item(68, 14).Value = -4614.363499999999
$ws.Cells.Item(71, 8).Value = 2723.75
$ws.Cells.Item(71, 9).Value = 1860
$ws.Cells.Item(71, 10).Value = 3116.3635
$ws.Cells.Item(71, 11).Value = 9300
$ws.Cells.Item(71, 12).Value = 15581.8175
$ws.Cells.Item(71, 13).Value = -5556
$ws.Cells.Item(71, 14).Value = -23069.8175
$ws.Cells.Item(132, 8).Value = 2941.8298
$ws.Cells.Item(132, 9).Value = 2248.9312
$ws.Cells.Item(132, 10).Value = 4058.1667
$ws.Cells.Item(132, 11).Value = 6746.7936
$ws.Cells.Item(132, 12).Value = 12174.5001
$ws.Cells.Item(132, 13).Value = -4216.7936
$ws.Cells.Item(132, 14).Value = -17234.5001
$ws.Cells.Item(134, 8).Value = 49964.777
$ws.Cells.Item(134, 10).Value = 49964.777
$ws.Cells.Item(134, 12).Value = 49964.777
$ws.Cells.Item(134, 14).Value = -60104.777
$ws.Cells.Item(137, 8).Value = 40783.332
$ws.Cells.Item(137, 10).Value = 40783.332
$ws.Cells.Item(137, 12).Value = 40783.332
$ws.Cells.Item(137, 14).Value = -50983.332
$ws.Cells.Item(139, 8).Value = 50199.6
$ws.Cells.Item(139, 10).Value = 50199.6
$ws.Cells.Item(139, 12).Value = 50199.6
$ws.Cells.Item(139, 14).Value = -60479.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 42892.5
$ws.Cells.Item(46, 10).Value = 42892.5
$ws.Cells.Item(46, 12).Value = 42892.5
$ws.Cells.Item(46, 14).Value = -43354.5
$ws.Cells.Item(62, 8).Value = 3229.0557
$ws.Cells.Item(62, 9).Value = 3166.5
$ws.Cells.Item(62, 10).Value = 3236.875
$ws.Cells.Item(62, 11).Value = 3166.5
$ws.Cells.Item(62, 12).Value = 3236.875
$ws.Cells.Item(62, 13).Value = -2542.5
$ws.Cells.Item(62, 14).Value = -4484.875
$ws.Cells.Item(65, 8).Value = 3229.0557
$ws.Cells.Item(65, 9).Value = 3166.5
$ws.Cells.Item(65, 10).Value = 3236.875
$ws.Cells.Item(65, 11).Value = 15832.5
$ws.Cells.Item(65, 12).Value = 16184.375
$ws.Cells.Item(65, 13).Value = -12712.5
$ws.Cells.Item(65, 14).Value = -22424.375
$ws.Cells.Item(134, 8).Value = 42892.5
$ws.Cells.Item(134, 10).Value = 42892.5
$ws.Cells.Item(134, 12).Value = 128677.5
$ws.Cells.Item(134, 14).Value = -133747.5
$ws.Cells.Item(139, 8).Value = 50999.8
$ws.Cells.Item(139, 10).Value = 50999.8
$ws.Cells.Item(139, 12).Value = 50999.8
$ws.Cells.Item(139, 14).Value = -61279.8
